# Updates cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.123.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.370.41"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.92%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("E9").Value = "  -2.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.35"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.29"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.37%  "

$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.738.37"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.367.13"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.096.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.66"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("E28").Value = "  +15.17%  "

$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.16"
$ws.Range("D30").ClearFormats()

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.64"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("E34").Value = "  +8.85%  "

$ws.Range("E36").Value = "  +1.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "127.45"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.06%  "

$ws.Range("E38").Value = "  +3.62%  "

$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("E42").Value = "  -5.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.929.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("E45").Value = "  +2.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.20"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -9.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.598.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("E49").Value = "  +2.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.43%  "

# Row 51: TrustWalletToken -> MultiversX
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.22%  "
